$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp banner (A1)
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 7 de Abril de 2020 a las 06:22'

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 2).Value = 367629
$ws.Cells.Item(4, 3).Value = 625
$ws.Cells.Item(4, 4).Value = 19810
$ws.Cells.Item(4, 5).Value = 336878
$ws.Cells.Item(4, 6).Value = 8983
$ws.Cells.Item(4, 7).Value = 70
$ws.Cells.Item(4, 8).Value = 10941

# Row 7: Alemania -> Alemania
$ws.Cells.Item(7, 4).Value = 36081
$ws.Cells.Item(7, 5).Value = 65484

# Row 29: India -> India
$ws.Cells.Item(29, 4).Value = 382
$ws.Cells.Item(29, 5).Value = 4260

# Row 108: Senegal -> Kirguistan
$ws.Cells.Item(108, 1).Value = 'Kirguistan'
$ws.Cells.Item(108, 2).Value = 228
$ws.Cells.Item(108, 3).Value = 12
$ws.Cells.Item(108, 4).Value = 33
$ws.Cells.Item(108, 5).Value = 191
$ws.Cells.Item(108, 6).Value = 5
$ws.Cells.Item(108, 8).Value = 4

# Row 109: Kirguistan -> Senegal
$ws.Cells.Item(109, 1).Value = 'Senegal'
$ws.Cells.Item(109, 2).Value = 226
$ws.Cells.Item(109, 4).Value = 92
$ws.Cells.Item(109, 5).Value = 132
$ws.Cells.Item(109, 6).Value = 1
$ws.Cells.Item(109, 8).Value = 2

# Row 125: Camboya -> Camboya
$ws.Cells.Item(125, 2).Value = 115
$ws.Cells.Item(125, 3).Value = 1
$ws.Cells.Item(125, 4).Value = 58
$ws.Cells.Item(125, 5).Value = 57

# Row 134: Guayana Francesa -> Guatemala
$ws.Cells.Item(134, 1).Value = 'Guatemala'
$ws.Cells.Item(134, 2).Value = 74
$ws.Cells.Item(134, 3).Value = 4
$ws.Cells.Item(134, 4).Value = 17
$ws.Cells.Item(134, 5).Value = 54
$ws.Cells.Item(134, 6).Value = 3
$ws.Cells.Item(134, 8).Value = 3

# Row 135: Aruba -> Guayana Francesa
$ws.Cells.Item(135, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(135, 2).Value = 72
$ws.Cells.Item(135, 4).Value = 34
$ws.Cells.Item(135, 5).Value = 38
$ws.Cells.Item(135, 6).Value = 1

# Row 136: Guatemala -> Aruba
$ws.Cells.Item(136, 1).Value = 'Aruba'
$ws.Cells.Item(136, 2).Value = 71
$ws.Cells.Item(136, 4).Value = 2
$ws.Cells.Item(136, 5).Value = 69
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 8).Value = 0

# Row 167: Guinea Ecuatorial -> Namibia
$ws.Cells.Item(167, 1).Value = 'Namibia'

# Row 168: Namibia -> Guinea Ecuatorial
$ws.Cells.Item(168, 1).Value = 'Guinea Ecuatorial'

# Row 170: Fiyi -> Antigua y Barbuda
$ws.Cells.Item(170, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 6).Value = 1

# Row 171: Antigua y Barbuda -> Fiyi
$ws.Cells.Item(171, 1).Value = 'Fiyi'
$ws.Cells.Item(171, 3).Value = 1
$ws.Cells.Item(171, 6).Value = 0

# Row 173: Mongolia -> Mongolia
$ws.Cells.Item(173, 4).Value = 4
$ws.Cells.Item(173, 5).Value = 11

# Row 183: Zimbabue -> Surinam
$ws.Cells.Item(183, 1).Value = 'Surinam'

# Row 184: Surinam -> Zimbabue
$ws.Cells.Item(184, 1).Value = 'Zimbabue'

# Row 194: Belice -> Somalia
$ws.Cells.Item(194, 1).Value = 'Somalia'
$ws.Cells.Item(194, 4).Value = 1
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 8).Value = 0

# Row 195: Somalia -> Belice
$ws.Cells.Item(195, 1).Value = 'Belice'
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 6).Value = 1
$ws.Cells.Item(195, 8).Value = 1

# Row 199: Botsuana -> San Bartolome
$ws.Cells.Item(199, 1).Value = 'San Bartolome'
$ws.Cells.Item(199, 4).Value = 1
$ws.Cells.Item(199, 8).Value = 0

# Row 200: San Bartolome -> Botsuana
$ws.Cells.Item(200, 1).Value = 'Botsuana'
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 8).Value = 1

# Row 207: Burundi -> Anguila
$ws.Cells.Item(207, 1).Value = 'Anguila'

# Row 208: Anguila -> Islas Virgenes Britanicas
$ws.Cells.Item(208, 1).Value = 'Islas Virgenes Britanicas'

# Row 209: Islas Virgenes Britanicas -> Burundi
$ws.Cells.Item(209, 1).Value = 'Burundi'

# Row 211: Bonaire, San Eustaquio y Saba -> Papua Nueva Guinea
$ws.Cells.Item(211, 1).Value = 'Papua Nueva Guinea'

# Row 212: Papua Nueva Guinea -> Bonaire, San Eustaquio y Saba
$ws.Cells.Item(212, 1).Value = 'Bonaire, San Eustaquio y Saba'

# Row 213: Timor Oriental -> Sudan del Sur
$ws.Cells.Item(213, 1).Value = 'Sudan del Sur'

# Row 214: Sudan del Sur -> San Pedro y Miquelon
$ws.Cells.Item(214, 1).Value = 'San Pedro y Miquelon'

# Row 215: San Pedro y Miquelon -> Timor Oriental
$ws.Cells.Item(215, 1).Value = 'Timor Oriental'
